$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 818.2941
$ws.Range("I80").Value = 874.4
$ws.Range("K80").Value = 2623.2
$ws.Range("M80").Value = -1625.2
$ws.Range("H83").Value = 818.2941
$ws.Range("I83").Value = 874.4
$ws.Range("K83").Value = 7869.599999999999
$ws.Range("M83").Value = -2877.599999999999
$ws.Range("H86").Value = 8336282
$ws.Range("I86").Value = 12502923
$ws.Range("J86").Value = 2999.5
$ws.Range("K86").Value = 12502923
$ws.Range("L86").Value = 2999.5
$ws.Range("M86").Value = -12501800
$ws.Range("N86").Value = -5245.5
$ws.Range("H89").Value = 8336282
$ws.Range("I89").Value = 12502923
$ws.Range("J89").Value = 2999.5
$ws.Range("K89").Value = 62514615
$ws.Range("L89").Value = 14997.5
$ws.Range("M89").Value = -62508999
$ws.Range("N89").Value = -26229.5
$ws.Range("H100").Value = 3542.7144
$ws.Range("I100").Value = 3143.5
$ws.Range("J100").Value = 4540.75
$ws.Range("K100").Value = 3143.5
$ws.Range("L100").Value = 4540.75
$ws.Range("M100").Value = -2602.5
$ws.Range("N100").Value = -5622.75
$ws.Range("H113").Value = 4914.5835
$ws.Range("J113").Value = 5000
$ws.Range("L113").Value = 5000
$ws.Range("N113").Value = -11508
$ws.Range("H125").Value = 1151.3846
$ws.Range("I125").Value = 947.5
$ws.Range("J125").Value = 1212.55
$ws.Range("K125").Value = 8527.5
$ws.Range("L125").Value = 10912.95
$ws.Range("M125").Value = -6067.5
$ws.Range("N125").Value = -15832.95
$ws.Range("H135").Value = 1358.7273
$ws.Range("I135").Value = 1242.2858
$ws.Range("K135").Value = 11180.5722
$ws.Range("M135").Value = -8645.572200000001
$ws.Range("H138").Value = 2825.5862
$ws.Range("J138").Value = 5106.5835
$ws.Range("L138").Value = 15319.7505
$ws.Range("N138").Value = -25599.7505
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5247.5
$ws.Range("J2").Value = 5247.5
$ws.Range("L2").Value = 5247.5
$ws.Range("N2").Value = -5473.5
$ws.Range("H34").Value = 30000
$ws.Range("J34").Value = 30000
$ws.Range("L34").Value = 30000
$ws.Range("N34").Value = -30542
$ws.Range("H61").Value = 3405.4546
$ws.Range("I61").Value = 3405.4546
$ws.Range("K61").Value = 3405.4546
$ws.Range("M61").Value = -3193.4546
$ws.Range("H116").Value = 5247.5
$ws.Range("J116").Value = 5247.5
$ws.Range("L116").Value = 5247.5
$ws.Range("N116").Value = -9835.5
$ws.Range("H132").Value = 1786.9231
$ws.Range("I132").Value = 1785.8334
$ws.Range("K132").Value = 5357.5002
$ws.Range("M132").Value = -2827.5002
$ws.Range("H136").Value = 3405.4546
$ws.Range("I136").Value = 3405.4546
$ws.Range("K136").Value = 10216.3638
$ws.Range("M136").Value = -7666.363799999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5247.5
$ws.Range("J3").Value = 5247.5
$ws.Range("L3").Value = 5247.5
$ws.Range("N3").Value = -5475.5
$ws.Range("H86").Value = 1717
$ws.Range("I86").Value = 1751.909
$ws.Range("J86").Value = 1640.2
$ws.Range("K86").Value = 1751.909
$ws.Range("L86").Value = 1640.2
$ws.Range("M86").Value = -628.9090000000001
$ws.Range("N86").Value = -3886.2
$ws.Range("H89").Value = 1717
$ws.Range("I89").Value = 1751.909
$ws.Range("J89").Value = 1640.2
$ws.Range("K89").Value = 8759.545
$ws.Range("L89").Value = 8201
$ws.Range("M89").Value = -3143.545
$ws.Range("N89").Value = -19433
$ws.Range("H107").Value = 9818.4
$ws.Range("I107").Value = 14364
$ws.Range("K107").Value = 14364
$ws.Range("M107").Value = -12444
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3062
$ws.Range("I16").Value = 1870.3334
$ws.Range("K16").Value = 1870.3334
$ws.Range("M16").Value = -1583.3334
$ws.Range("H103").Value = 21728.6
$ws.Range("I103").Value = 21910.75
$ws.Range("K103").Value = 21910.75
$ws.Range("M103").Value = -20738.75
$ws.Range("H105").Value = 3516.6562
$ws.Range("I105").Value = 3845
$ws.Range("K105").Value = 3845
$ws.Range("M105").Value = -2098
$ws.Range("H107").Value = 864.2143
$ws.Range("I107").Value = 430.125
$ws.Range("J107").Value = 1443
$ws.Range("K107").Value = 430.125
$ws.Range("L107").Value = 1443
$ws.Range("M107").Value = 1489.875
$ws.Range("N107").Value = -5283
$ws.Range("H113").Value = 3062
$ws.Range("I113").Value = 1870.3334
$ws.Range("K113").Value = 1870.3334
$ws.Range("M113").Value = 299.6666
$ws.Range("H134").Value = 2205.25
$ws.Range("I134").Value = 1794.9
$ws.Range("K134").Value = 5384.700000000001
$ws.Range("M134").Value = -2849.700000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 11595.8
$ws.Range("I62").Value = 9333
$ws.Range("J62").Value = 14990
$ws.Range("K62").Value = 27999
$ws.Range("L62").Value = 44970
$ws.Range("M62").Value = -27313
$ws.Range("N62").Value = -46342
$ws.Range("H65").Value = 11595.8
$ws.Range("I65").Value = 9333
$ws.Range("J65").Value = 14990
$ws.Range("K65").Value = 83997
$ws.Range("L65").Value = 134910
$ws.Range("M65").Value = -80565
$ws.Range("N65").Value = -141774
$ws.Range("H138").Value = 3499.8
$ws.Range("I138").Value = 3374.75
$ws.Range("J138").Value = 4000
$ws.Range("K138").Value = 10124.25
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = -4984.25
$ws.Range("N138").Value = -22280
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9444.5625
$ws.Range("I80").Value = 9932.615
$ws.Range("K80").Value = 9932.615
$ws.Range("M80").Value = -8934.615
$ws.Range("H83").Value = 9444.5625
$ws.Range("I83").Value = 9932.615
$ws.Range("K83").Value = 49663.075
$ws.Range("M83").Value = -44671.075
$ws.Range("H107").Value = 27836.166
$ws.Range("I107").Value = 34488.57
$ws.Range("J107").Value = 4552.75
$ws.Range("K107").Value = 34488.57
$ws.Range("L107").Value = 4552.75
$ws.Range("M107").Value = -32568.57
$ws.Range("N107").Value = -8392.75
$ws.Range("H113").Value = 125003624
$ws.Range("I113").Value = 166669500
$ws.Range("K113").Value = 166669500
$ws.Range("M113").Value = -166667330
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 73103.336
$ws.Range("I46").Value = 86924
$ws.Range("K46").Value = 86924
$ws.Range("M46").Value = -86736
$ws.Range("H132").Value = 2399.389
$ws.Range("I132").Value = 1824.625
$ws.Range("J132").Value = 6997.5
$ws.Range("K132").Value = 5473.875
$ws.Range("L132").Value = 20992.5
$ws.Range("M132").Value = -2943.875
$ws.Range("N132").Value = -26052.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 147451.55
$ws.Range("I62").Value = 5493.25
$ws.Range("K62").Value = 5493.25
$ws.Range("M62").Value = -4869.25
$ws.Range("H65").Value = 147451.55
$ws.Range("I65").Value = 5493.25
$ws.Range("K65").Value = 27466.25
$ws.Range("M65").Value = -24346.25
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H126").Value = 2477.1155
$ws.Range("I126").Value = 2293.35
$ws.Range("J126").Value = 3089.6667
$ws.Range("K126").Value = 6880.049999999999
$ws.Range("L126").Value = 9269.000100000001
$ws.Range("M126").Value = -4410.049999999999
$ws.Range("N126").Value = -14209.0001
